$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Edit 1 (Road network section): "RoadNetwork_0 " -> "RoadNetwork_1"
#   - the run holding the trailing space character is removed entirely
#   - the run holding "0" keeps its own identity, its text becomes "1"
# ------------------------------------------------------------------
$f1 = $d.Content
$f1.Find.Execute("RoadNetwork_0 ", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$matchEnd = $f1.End

$spaceRange = $d.Range($matchEnd - 1, $matchEnd)
$spaceRange.Delete()

$zeroRange = $d.Range($matchEnd - 2, $matchEnd - 1)
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidR="00406D2A"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:lang w:val="vi-VN"/></w:rPr><w:t>1</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$zeroRange.InsertXML($xml1)

# ------------------------------------------------------------------
# Edit 2 (Land use section): " LandUse_0" -> " LandUse_" + new run "1"
# ------------------------------------------------------------------
$f2 = $d.Content
$f2.Find.Execute(" LandUse_0", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$rng2 = $d.Range($f2.Start, $f2.End)
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:lang w:val="vi-VN"/></w:rPr><w:t xml:space="preserve"> LandUse_</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:lang w:val="vi-VN"/></w:rPr><w:t>1</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng2.InsertXML($xml2)

# ------------------------------------------------------------------
# Edit 3 (Building shapefile name): "Building_0.shp" ->
#   "Building_" (same run) + new run "1" + new run ".shp"
# ------------------------------------------------------------------
$f3 = $d.Content
$f3.Find.Execute("Building_0.shp", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$rng3 = $d.Range($f3.Start, $f3.End)
$xml3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidRPr="004E2175"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="000000"/><w:lang w:val="vi-VN"/></w:rPr><w:t>Building_</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="000000"/><w:lang w:val="vi-VN"/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="000000"/><w:lang w:val="vi-VN"/></w:rPr><w:t>.shp</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng3.InsertXML($xml3)

Write-Output "Scenario 1 shapefile names updated (RoadNetwork_1, LandUse_1, Building_1.shp)"
